# Slide 13: update the two phone/number placeholder shapes.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)

# --- "Rectangle 28" (shape 6): "Text 1267920" -> "Text 11888" ---
$sh1 = $s.Shapes.Item(6)
$tr1 = $sh1.TextFrame.TextRange

# Split the trailing space off the non-bold "Text " run into its own run
# (re-assert the same font name to force a run break without changing formatting).
$space1 = $tr1.Characters(5, 1)
$space1.Font.Name = "Avenir Light"

# Replace the bold number run's text.
$num1 = $tr1.Characters(6, 7)
$num1.Text = "11888"

# --- "Rectangle 29" (shape 7): "to (760) 452 -8549" -> "to (760) 452 -8548" ---
$sh2 = $s.Shapes.Item(7)
$tr2 = $sh2.TextFrame.TextRange

# Split the last 4 digits of the bold run into their own run, then replace the text.
$tail2 = $tr2.Characters(15, 4)
$tail2.Font.Name = "Avenir Light"
$tail2.Text = "8548"
